$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOQ4023")

# Semestre ideal: EQD-8,EQN-9 -> EQD-9,EQN-10
$ws.Range("B9").Value = "EQD-9,EQN-10"
$ws.Range("C9").Value = "EQD-9,EQN-10"

# Docentes responsáveis: 5840855 - Heizir Ferreira de Castro -> 1285870 - Marcos Villela Barcza
$ws.Range("B13").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C13").Value = "1285870 - Marcos Villela Barcza"

# Requisitos (first): LOQ4002 -  Reatores Quimicos  (Requisito fraco) -> LOQ4038 -  Química Orgânica II  (Requisito fraco)
$req1 = "LOQ4038 -  Química Orgânica II  (Requisito fraco)`n"
$ws.Range("B24").Value = $req1
$ws.Range("C24").Value = $req1

# Requisitos (second): LOT2004 -  Bioquímica  (Requisito fraco) -> LOQ4057 -  Operações Unitárias III  (Requisito fraco)
$req2 = "LOQ4057 -  Operações Unitárias III  (Requisito fraco)`n"
$ws.Range("B25").Value = $req2
$ws.Range("C25").Value = $req2
